$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.497.22"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.843.76"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  +0.16%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "261.66"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.5336"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +2.34%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3039"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -5.90%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06895"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +1.57%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "18.13"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -2.87%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07681"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.867.80"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.7402"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -4.47%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "89.79"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +1.73%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.997"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("E16").Value = "  +0.19%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "13.98"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000007950"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "26.518.45"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "2.077.16"
$ws.Range("E21").Value = "  -0.51%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.618"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.992"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "9.308"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.42%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "143.29"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.32%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.201"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.83%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "1.687"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.35%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "16.99"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.03%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "110.80"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.89%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "4.262"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +2.50%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.08790"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.69%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.064"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.94%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.04802"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "2.932"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +2.48%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.7274"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("E36").Value = "  +0.82%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "3.110"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.68%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.306"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +4.75%  "
$ws.Range("E39").Value = "  -4.08%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.4774"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.28%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.9158"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +3.21%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "108.02"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.82%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "5.881"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -2.58%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.17%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "7.484"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.62%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "9.090"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.4125"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.22%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.1242"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.72%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "34.84"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.05799"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.8970"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.02%  "
